# Generate Report for Handoff
#
# Updates the "Latest Handoff Date/Datetime" for the file
# 1ebd9f8b-1134-4c3c-ab15-1c28586f6d2f (status "Ready for handoff") to the
# freshly generated handoff timestamps, on the Overview sheet and on each
# per-locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-30-19 06:30:51"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-19 06:30:49"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-19 06:30:51"
